# "Add files via upload" — the author re-uploaded an updated copy of the
# "Chiffres COVID-19 Valais" tracker. The substantive change is new/revised
# daily figures for rows 542-559 (columns C, E, F, G, L, M — the raw inputs).
# Columns B, H, J and K are formulas (cumulative cases, SI+hosp total,
# cumulative deaths, new deaths) and recompute automatically once the
# inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 542 — 20.08.2021: new positive cases 91 -> 92
$ws.Range("C542").Value = 92

# Row 543 — 21.08.2021: new positive cases 61 -> 62
$ws.Range("C543").Value = 62

# Row 549 — 27.08.2021: new positive cases 81 -> 80
$ws.Range("C549").Value = 80

# Row 552 — 30.08.2021: new positive cases 121 -> 124, SI patients 19 -> 18
$ws.Range("C552").Value = 124
$ws.Range("G552").Value = 18

# Row 553 — 31.08.2021: new positive cases 71 -> 73, new admissions 10 -> 11,
# SI patients 19 -> 16
$ws.Range("C553").Value = 73
$ws.Range("E553").Value = 11
$ws.Range("G553").Value = 16

# Row 554 — 01.09.2021: new positive cases 67 -> 90, new admissions 11 -> 10,
# SI patients 15 -> 16
$ws.Range("C554").Value = 90
$ws.Range("E554").Value = 10
$ws.Range("G554").Value = 16

# Row 555 — 02.09.2021: new positive cases 3 -> 108, new admissions 10 -> 9,
# SI patients 14 -> 15
$ws.Range("C555").Value = 108
$ws.Range("E555").Value = 9
$ws.Range("G555").Value = 15

# Row 556 — 03.09.2021: previously blank, now filled in with a full day of data
$ws.Range("C556").Value = 91
$ws.Range("E556").Value = 9
$ws.Range("F556").Value = 6
$ws.Range("G556").Value = 16

# Row 557 — 04.09.2021: previously blank, now filled in
$ws.Range("C557").Value = 53
$ws.Range("E557").Value = 9
$ws.Range("F557").Value = 7
$ws.Range("G557").Value = 13

# Row 558 — 05.09.2021: previously blank, now filled in
$ws.Range("C558").Value = 31
$ws.Range("E558").Value = 10
$ws.Range("F558").Value = 8
$ws.Range("G558").Value = 13

# Row 559 — 06.09.2021: previously blank, now filled in
$ws.Range("C559").Value = 9
$ws.Range("E559").Value = 10
$ws.Range("F559").Value = 8
$ws.Range("G559").Value = 15

# Columns L ("Nb nouveaux décès à l'hôpital") and M ("... extra-hospitaliers")
# are formatted as Text (@), so a plain .Value assignment of 0 would be
# stored as the string "0" instead of the number 0. Briefly switch those
# cells to a numeric format, write the number, then restore the original
# Text format so the cell keeps its original style/appearance.
foreach ($addr in @("L556","M556","L557","M557","L558","M558","L559","M559")) {
    $cell = $ws.Range($addr)
    $savedFormat = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = 0
    $cell.NumberFormat = $savedFormat
}
